# Update automation results for the DNF system (subj0 / results_Alpha)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R2").Value = 28.44618809978375
$ws.Range("R3").Value = 26.0805135663567

$ws.Range("R4").Value = 28.08518718029109
$ws.Range("S4").Value = 18.76295454903327

$ws.Range("R5").Value = 1.553820082092557
$ws.Range("S5").Value = -0.4416577113042988

$ws.Range("R6").Value = 25.64080623124418
$ws.Range("S6").Value = -9.408754405736811
